$d = $word.ActiveDocument

# Update the date heading paragraph
$d.Content.Find.Execute("2024-09-09 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-09-10 Tuesday", 2)

# Update the 100 arithmetic-expression cells (row-major order, 20 rows x 5 cols)
$newValues = @"
82-36=
49+38=
18+38=
51-45=
59+18=
96-7=
86-19=
40-28=
19+38=
9+73=
52-35=
19+23=
59+14=
34+9=
36+25=
71-7=
45+37=
93-88=
29+46=
47+17=
27+36=
73-26=
93-56=
95-38=
85+9=
55+16=
41-25=
23+48=
92-35=
7+45=
88-69=
9+52=
46+39=
83-28=
95-79=
73-15=
61-52=
53-7=
17+17=
57+39=
73+18=
88-59=
90-83=
57-28=
61-53=
56+8=
90-73=
59+14=
94-28=
41-5=
73-64=
42-35=
35-27=
81-43=
70-19=
93-24=
62+29=
88-69=
15-7=
62-38=
59+7=
55-38=
82-36=
76-68=
17+76=
17+47=
4+68=
17+65=
39+3=
92-79=
29+46=
94-69=
62-54=
96-47=
59+3=
67+26=
17+34=
38+54=
58+19=
80-41=
87-68=
4+49=
86+6=
8+66=
84-19=
5+68=
8+63=
57-49=
81-13=
55+38=
8+39=
37+35=
83-26=
81-77=
58+19=
7+24=
65-58=
53-29=
17+75=
18+17=
"@ -split "`r?`n" | Where-Object { $_ -ne "" }

$t = $d.Tables.Item(1)

if ($newValues.Count -ne 100) {
    Write-Output "WARNING: expected 100 replacement values, found $($newValues.Count)"
}
if (($t.Rows.Count * $t.Columns.Count) -ne 100) {
    Write-Output "WARNING: expected a 20x5 (100-cell) table, found $($t.Rows.Count)x$($t.Columns.Count)"
}

$idx = 0
foreach ($row in $t.Rows) {
    foreach ($cell in $row.Cells) {
        $cell.Range.Text = $newValues[$idx]
        $idx = $idx + 1
    }
}
Write-Output "updated cells: $idx"
